# Add "joelito/legal-xlm-longformer-base" as a new model row at the
# bottom of the list (row 41), matching the existing layout used by all
# the other model rows:
#   A = _name_or_path
#   B = revision
#   C = need_to_be_run_with_LEXTREME
#   D = need_to_be_run_with_LexGlue

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 41

$ws.Cells.Item($newRow, 1).Value = "joelito/legal-xlm-longformer-base"
$ws.Cells.Item($newRow, 2).Value = "main"
$ws.Cells.Item($newRow, 3).Value = $true
$ws.Cells.Item($newRow, 4).Value = $true

# Leave the view scrolled to / focused on the newly added row, the way
# the author's session ended up (selection C41:D41).
$ws.Range("C41:D41").Select()
